$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2021-10-11"

# Update row 12 label
$ws.Range("A12").Value = "October (through 10-11)"

# Update row 12 values
$ws.Range("C12").Value = 9
$ws.Range("D12").Value = 0.1
$ws.Range("F12").Value = 18
$ws.Range("I12").Value = 14
$ws.Range("J12").Value = 0.2222
$ws.Range("L12").Value = 27
$ws.Range("M12").Value = 0.06900000000000001
$ws.Range("R12").Value = 50
$ws.Range("U12").Value = 72

# Update row 13 (Total) values
$ws.Range("C13").Value = 205
$ws.Range("D13").Value = 0.1314
$ws.Range("F13").Value = 401
$ws.Range("G13").Value = 0.1029
$ws.Range("I13").Value = 591
$ws.Range("J13").Value = 0.0837
$ws.Range("L13").Value = 514
$ws.Range("M13").Value = 0.1092
$ws.Range("R13").Value = 898
$ws.Range("S13").Value = 0.0557
$ws.Range("U13").Value = 1243
$ws.Range("V13").Value = 0.0598
